$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6

# Row 3
$ws.Range("G3").Value = 2.88
$ws.Range("I3").Value = 2.88
$ws.Range("O3").Value = 1.73
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 3.5
$ws.Range("R3").Value = 1.3
$ws.Range("S3").Value = 1.78
$ws.Range("T3").Value = 2.03
$ws.Range("AA3").Value = 34
$ws.Range("AI3").Value = 12

# Row 4
$ws.Range("G4").Value = 2.2
$ws.Range("H4").Value = 2.88
$ws.Range("I4").Value = 4.1
$ws.Range("J4").Value = 3.1
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = 1.73
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 3.5
$ws.Range("R4").Value = 1.3
$ws.Range("S4").Value = 1.75
$ws.Range("T4").Value = 2.05
$ws.Range("X4").Value = 8.5
$ws.Range("Z4").Value = 21
$ws.Range("AA4").Value = 26
$ws.Range("AC4").Value = 4.75
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 23
$ws.Range("AF4").Value = 101
$ws.Range("AJ4").Value = 17
$ws.Range("AM4").Value = 67
$ws.Range("AO4").Value = 15
$ws.Range("AS4").Value = 451
$ws.Range("AT4").Value = 2
$ws.Range("BA4").Value = 201

# Row 6
$ws.Range("S6").Value = 1.47

# Row 9
$ws.Range("G9").Value = 3.25
$ws.Range("H9").Value = 3.1
$ws.Range("I9").Value = 2.3
$ws.Range("L9").Value = 3
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.62
$ws.Range("AA9").Value = 29
$ws.Range("AH9").Value = 7
$ws.Range("AJ9").Value = 9.5
$ws.Range("AK9").Value = 21
$ws.Range("AO9").Value = 19
$ws.Range("AY9").Value = 23

# Row 10
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 6

# Row 11
$ws.Range("G11").Value = 1.65
$ws.Range("I11").Value = 5.75
$ws.Range("L11").Value = 6
$ws.Range("N11").Value = 7.5
$ws.Range("O11").Value = 1.4
$ws.Range("P11").Value = 2.75
$ws.Range("U11").Value = 2.25
$ws.Range("V11").Value = 1.57
$ws.Range("W11").Value = 5.5
$ws.Range("Z11").Value = 12
$ws.Range("AC11").Value = 7.5
$ws.Range("AE11").Value = 21
$ws.Range("AH11").Value = 12
$ws.Range("AJ11").Value = 19
$ws.Range("AO11").Value = 9
$ws.Range("AP11").Value = 23
$ws.Range("AS11").Value = 201
